# === GRAL.xlsx update: various additions, clean up of older models ===

$wb = $excel.ActiveWorkbook
$main  = $wb.Worksheets.Item("Main")
$model = $wb.Worksheets.Item("Model")

# -----------------------------------------------------------------
# Main sheet ("Main")
# -----------------------------------------------------------------

# Header cell moved from C2 to B2
$main.Range("C2").Copy($main.Range("B2"))
$main.Range("C2").Clear()

# New right-aligned "Q224" labels next to the two totals
$main.Range("L9").Value = "Q224"
$main.Range("L9").HorizontalAlignment = -4152
$main.Range("L10").Value = "Q224"
$main.Range("L10").HorizontalAlignment = -4152

# View changes: zoom + selection (done without leaving Model as the active sheet)
$main.Activate()
$excel.ActiveWindow.Zoom = 145
$main.Range("J8").Select()
$model.Activate()

# -----------------------------------------------------------------
# Model sheet ("Model")
# -----------------------------------------------------------------

# View changes: zoom + bottom-right pane selection
$model.Activate()
$excel.ActiveWindow.Zoom = 160
$model.Range("H17").Select()

# Row 3 ("Screening") promoted to the bold subtotal-style formatting
# already used by row 6 ("Revenue")
$model.Range("B3:J3").Copy()
$model.Range("B6:J6").PasteSpecial(-4122)

# Convert the year header fill-across (N2:Z2) into one relative-fill formula
$model.Range("N2:Z2").Formula = "=+M2+1"

# H19: baseline cash balance now derives from a small formula instead of a literal
$model.Range("H19").Formula = "=959+3.918"
$model.Range("H19").NumberFormat = "#,##0"

# Turn the previously-independent O8:T8 / O13:T13 formula rows into shared fills
$model.Range("O8:T8").Formula = "=+O6-O7"
$model.Range("O13:T13").Formula = "=+O8-O12"

# --- New balance-sheet style rows 20-27 (order chosen to control shared-string ids) ---
$model.Range("B27").Value = "Assets"
$model.Range("B26").Value = "Goodwill"
$model.Range("B20").Value = "AR"
$model.Range("B21").Value = "Supplies"
$model.Range("B22").Value = "Prepaids"
$model.Range("B23").Value = "PP&E"
$model.Range("B24").Value = "Lease"
$model.Range("B25").Value = "ONCA"

$model.Range("H20").Formula = "=13.374+0.032"
$model.Range("H21").Formula = "=18.196+7.31"
$model.Range("H22").Formula = "=20.866+0.059"
$model.Range("H23").Formula = "=74.984+3.021"
$model.Range("H24").Value = 74.503
$model.Range("H25").Value = 8.476
$model.Range("H26").Value = 2086.056
$model.Range("H27").Formula = "=SUM(H19:H26)"

$model.Range("H20:H27").NumberFormat = "#,##0"
$model.Range("H20:H27").HorizontalAlignment = -4152

# --- New standalone value rows 29-35 ---
$model.Range("H29").Value = 16.247
$model.Range("H30").Value = 56.5473
$model.Range("H31").Value = 13.945
$model.Range("H32").Value = 1.413
$model.Range("H33").Value = 62.165
$model.Range("H34").Value = 422.163
$model.Range("H35").Value = 2.007

$model.Range("H29:H35").NumberFormat = "#,##0"

Write-Output "edit applied"
